# Update Premier League team statistics with Betting Markets Analytics data refresh
# (Arsenal row 2 and Brentford row 8 advance by one additional match played)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 26
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 18
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 381
$ws.Range("N2").Value = 38
$ws.Range("P2").Value = 288
$ws.Range("Q2").Value = 93
$ws.Range("R2").Value = 11
$ws.Range("U2").Value = 76
$ws.Range("V2").Value = 54
$ws.Range("W2").Value = 52
$ws.Range("X2").Value = 129
$ws.Range("Y2").Value = 141
$ws.Range("Z2").Value = 111
$ws.Range("AA2").Value = 189
$ws.Range("AB2").Value = 388
$ws.Range("AC2").Value = 157
$ws.Range("AH2").Value = 57.461538461538
$ws.Range("AI2").Value = 12511
$ws.Range("AJ2").Value = 10625
$ws.Range("AK2").Value = 84.925265766126
$ws.Range("AL2").Value = 5427
$ws.Range("AM2").Value = 5015
$ws.Range("AN2").Value = 92.408328726737
$ws.Range("AO2").Value = 7084
$ws.Range("AP2").Value = 5610
$ws.Range("AQ2").Value = 79.192546583851
$ws.Range("AR2").Value = 758
$ws.Range("AS2").Value = 350
$ws.Range("AT2").Value = 46.174142480211
$ws.Range("AU2").Value = 503
$ws.Range("AV2").Value = 124
$ws.Range("AW2").Value = 24.652087475149
$ws.Range("AY2").Value = 410
$ws.Range("AZ2").Value = 184
$ws.Range("BA2").Value = 37
$ws.Range("BC2").Value = 14
$ws.Range("BF2").Value = 606
$ws.Range("BH2").Value = 8
$ws.Range("BI2").Value = 2447
$ws.Range("BJ2").Value = 1274
$ws.Range("BK2").Value = 52.063751532489
$ws.Range("BL2").Value = 1642
$ws.Range("BM2").Value = 859
$ws.Range("BN2").Value = 52.31425091352
$ws.Range("BO2").Value = 805
$ws.Range("BP2").Value = 415
$ws.Range("BQ2").Value = 51.552795031056
$ws.Range("BR2").Value = 3196
$ws.Range("BT2").Value = 263
$ws.Range("BU2").Value = 36
$ws.Range("BX2").Value = 6.9389312977099
$ws.Range("BY2").Value = 1465
$ws.Range("BZ2").Value = 3307
$ws.Range("CA2").Value = 4061
$ws.Range("CB2").Value = 7368
$ws.Range("CC2").Value = 29
$ws.Range("CD2").Value = 17
$ws.Range("CE2").Value = 19
$ws.Range("CF2").Value = 843
$ws.Range("CG2").Value = 81
$ws.Range("CH2").Value = 55
$ws.Range("CI2").Value = 321
$ws.Range("CJ2").Value = 333
$ws.Range("CK2").Value = 149
$ws.Range("CM2").Value = 25
$ws.Range("CO2").Value = 188
$ws.Range("CP2").Value = 148
$ws.Range("CQ2").Value = 550
$ws.Range("CR2").Value = 1393
$ws.Range("CU2").Value = 202
$ws.Range("CV2").Value = 71
$ws.Range("CW2").Value = 131
$ws.Range("CX2").Value = 71
$ws.Range("CY2").Value = 75
$ws.Range("CZ2").Value = 56
$ws.Range("DA2").Value = 71
$ws.Range("DB2").Value = 391
$ws.Range("DC2").Value = 2309
$ws.Range("DD2").Value = 4684
$ws.Range("DE2").Value = 4571
$ws.Range("DF2").Value = 9255
$ws.Range("DG2").Value = 45
$ws.Range("DH2").Value = 442
$ws.Range("DI2").Value = 128
$ws.Range("DJ2").Value = 1170
$ws.Range("DK2").Value = 266
$ws.Range("DM2").Value = 26
$ws.Range("DQ2").Value = 1.92
$ws.Range("DR2").Value = 0.69
$ws.Range("DS2").Value = 14.65
$ws.Range("DT2").Value = 4.96
$ws.Range("DU2").Value = 6.04
$ws.Range("DV2").Value = 10.12
$ws.Range("DW2").Value = 1.38
$ws.Range("DY2").Value = 2.92
$ws.Range("DZ2").Value = 2
$ws.Range("EA2").Value = 0.31
$ws.Range("EC2").Value = 408.65
$ws.Range("ED2").Value = 15.77
$ws.Range("EE2").Value = 7.08
$ws.Range("EF2").Value = 23.31
$ws.Range("EG2").Value = 1.42
$ws.Range("EH2").Value = 49
$ws.Range("EI2").Value = 33.04
$ws.Range("EJ2").Value = 15.96
$ws.Range("EK2").Value = 122.92
$ws.Range("EL2").Value = 7.27

# Row 8
$ws.Range("D8").Value = 26
$ws.Range("E8").Value = 40
$ws.Range("F8").Value = 35
$ws.Range("H8").Value = 24
$ws.Range("I8").Value = 270
$ws.Range("N8").Value = 38
$ws.Range("P8").Value = 206
$ws.Range("Q8").Value = 64
$ws.Range("R8").Value = 7
$ws.Range("U8").Value = 79
$ws.Range("V8").Value = 59
$ws.Range("W8").Value = 44
$ws.Range("X8").Value = 105
$ws.Range("Y8").Value = 100
$ws.Range("Z8").Value = 65
$ws.Range("AA8").Value = 160
$ws.Range("AB8").Value = 402
$ws.Range("AC8").Value = 128
$ws.Range("AE8").Value = 51
$ws.Range("AG8").Value = 43
$ws.Range("AH8").Value = 46.423076923077
$ws.Range("AI8").Value = 10057
$ws.Range("AJ8").Value = 7985
$ws.Range("AK8").Value = 79.397434622651
$ws.Range("AL8").Value = 4987
$ws.Range("AM8").Value = 4508
$ws.Range("AN8").Value = 90.395027070383
$ws.Range("AO8").Value = 5070
$ws.Range("AP8").Value = 3477
$ws.Range("AQ8").Value = 68.579881656805
$ws.Range("AR8").Value = 958
$ws.Range("AS8").Value = 419
$ws.Range("AT8").Value = 43.736951983299
$ws.Range("AU8").Value = 456
$ws.Range("AV8").Value = 118
$ws.Range("AW8").Value = 25.877192982456
$ws.Range("AY8").Value = 387
$ws.Range("AZ8").Value = 237
$ws.Range("BA8").Value = 69
$ws.Range("BC8").Value = 10
$ws.Range("BF8").Value = 777
$ws.Range("BI8").Value = 2636
$ws.Range("BJ8").Value = 1266
$ws.Range("BK8").Value = 48.027314112291
$ws.Range("BL8").Value = 1673
$ws.Range("BM8").Value = 775
$ws.Range("BN8").Value = 46.323968918111
$ws.Range("BO8").Value = 963
$ws.Range("BP8").Value = 491
$ws.Range("BQ8").Value = 50.986500519211
$ws.Range("BR8").Value = 3521
$ws.Range("BT8").Value = 280
$ws.Range("BU8").Value = 52
$ws.Range("BX8").Value = 6.801038961039
$ws.Range("BY8").Value = 2546
$ws.Range("BZ8").Value = 5063
$ws.Range("CA8").Value = 4566
$ws.Range("CB8").Value = 9629
$ws.Range("CC8").Value = 55
$ws.Range("CD8").Value = 36
$ws.Range("CE8").Value = 34
$ws.Range("CF8").Value = 840
$ws.Range("CG8").Value = 134
$ws.Range("CH8").Value = 99
$ws.Range("CI8").Value = 500
$ws.Range("CJ8").Value = 362
$ws.Range("CK8").Value = 183
$ws.Range("CM8").Value = 18
$ws.Range("CO8").Value = 186
$ws.Range("CP8").Value = 225
$ws.Range("CQ8").Value = 567
$ws.Range("CR8").Value = 1324
$ws.Range("CU8").Value = 328
$ws.Range("CV8").Value = 117
$ws.Range("CW8").Value = 208
$ws.Range("CX8").Value = 120
$ws.Range("CY8").Value = 108
$ws.Range("CZ8").Value = 103
$ws.Range("DA8").Value = 117
$ws.Range("DB8").Value = 453
$ws.Range("DC8").Value = 3541
$ws.Range("DD8").Value = 6616
$ws.Range("DE8").Value = 5104
$ws.Range("DF8").Value = 11720
$ws.Range("DG8").Value = 44
$ws.Range("DH8").Value = 502
$ws.Range("DI8").Value = 185
$ws.Range("DJ8").Value = 1302
$ws.Range("DK8").Value = 236
$ws.Range("DM8").Value = 26
$ws.Range("DQ8").Value = 1.54
$ws.Range("DR8").Value = 1.35
$ws.Range("DS8").Value = 10.38
$ws.Range("DT8").Value = 4.04
$ws.Range("DU8").Value = 4.92
$ws.Range("DV8").Value = 10.77
$ws.Range("DW8").Value = 2
$ws.Range("DZ8").Value = 1.69
$ws.Range("EA8").Value = 0.38
$ws.Range("EB8").Value = 0.23
$ws.Range("EC8").Value = 307.12
$ws.Range("ED8").Value = 14.88
$ws.Range("EE8").Value = 9.12
$ws.Range("EF8").Value = 29.88
$ws.Range("EG8").Value = 2.65
$ws.Range("EH8").Value = 48.69
$ws.Range("EI8").Value = 29.81
$ws.Range("EJ8").Value = 18.88
$ws.Range("EK8").Value = 135.42
$ws.Range("EL8").Value = 6.15

Write-Host "Updated $($wb.ActiveSheet.Name): 219 cells across rows 2 and 8"
